$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2537.5
$ws.Range("I2").Value = 1750
$ws.Range("J2").Value = 4900
$ws.Range("K2").Value = 1750
$ws.Range("L2").Value = 4900
$ws.Range("M2").Value = -1637
$ws.Range("N2").Value = -5126

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6055.7617
$ws.Range("I137").Value = 1098.9166
$ws.Range("J137").Value = 12664.889
$ws.Range("K137").Value = 3296.7498
$ws.Range("L137").Value = 37994.667
$ws.Range("M137").Value = -746.7498000000001
$ws.Range("N137").Value = -43094.667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2538.3823
$ws.Range("J141").Value = 3325.4167
$ws.Range("L141").Value = 9976.250100000001
$ws.Range("N141").Value = -20336.2501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2443.5833
$ws.Range("I61").Value = 912.3
$ws.Range("J61").Value = 10100
$ws.Range("K61").Value = 912.3
$ws.Range("L61").Value = 10100
$ws.Range("M61").Value = -700.3
$ws.Range("N61").Value = -10524

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1391.2162
$ws.Range("I74").Value = 1648.3889
$ws.Range("J74").Value = 1147.579
$ws.Range("K74").Value = 1648.3889
$ws.Range("L74").Value = 1147.579
$ws.Range("M74").Value = -774.3888999999999
$ws.Range("N74").Value = -2895.579

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1391.2162
$ws.Range("I77").Value = 1648.3889
$ws.Range("J77").Value = 1147.579
$ws.Range("K77").Value = 8241.9445
$ws.Range("L77").Value = 5737.895
$ws.Range("M77").Value = -3873.9445
$ws.Range("N77").Value = -14473.895

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 31720.838
$ws.Range("I132").Value = 47204.227
$ws.Range("J132").Value = 9011.866
$ws.Range("K132").Value = 141612.681
$ws.Range("L132").Value = 27035.598
$ws.Range("M132").Value = -139082.681
$ws.Range("N132").Value = -32095.598

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2443.5833
$ws.Range("I136").Value = 912.3
$ws.Range("J136").Value = 10100
$ws.Range("K136").Value = 2736.9
$ws.Range("L136").Value = 30300
$ws.Range("M136").Value = -186.8999999999996
$ws.Range("N136").Value = -35400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3546.3635
$ws.Range("I134").Value = 2259.2
$ws.Range("J134").Value = 4619
$ws.Range("K134").Value = 6777.599999999999
$ws.Range("L134").Value = 13857
$ws.Range("M134").Value = -4242.599999999999
$ws.Range("N134").Value = -18927

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21089.492
$ws.Range("I31").Value = 21877.14
$ws.Range("J31").Value = 16599.9
$ws.Range("K31").Value = 21877.14
$ws.Range("L31").Value = 16599.9
$ws.Range("M31").Value = -21582.14
$ws.Range("N31").Value = -17189.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 21089.492
$ws.Range("I34").Value = 21877.14
$ws.Range("J34").Value = 16599.9
$ws.Range("K34").Value = 21877.14
$ws.Range("L34").Value = 16599.9
$ws.Range("M34").Value = -21675.14
$ws.Range("N34").Value = -17003.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5045.3335
$ws.Range("I58").Value = 1417.7142
$ws.Range("J58").Value = 8219.5
$ws.Range("K58").Value = 1417.7142
$ws.Range("L58").Value = 8219.5
$ws.Range("M58").Value = -1214.7142
$ws.Range("N58").Value = -8625.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3299.6
$ws.Range("I132").Value = 2449.6667
$ws.Range("J132").Value = 4574.5
$ws.Range("K132").Value = 7349.000100000001
$ws.Range("L132").Value = 13723.5
$ws.Range("M132").Value = -4819.000100000001
$ws.Range("N132").Value = -18783.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 33335712
$ws.Range("I134").Value = 2074.889
$ws.Range("J134").Value = 83336170
$ws.Range("K134").Value = 6224.667
$ws.Range("L134").Value = 250008510
$ws.Range("M134").Value = -3689.667
$ws.Range("N134").Value = -250013580

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5045.3335
$ws.Range("I136").Value = 1417.7142
$ws.Range("J136").Value = 8219.5
$ws.Range("K136").Value = 4253.142599999999
$ws.Range("L136").Value = 24658.5
$ws.Range("M136").Value = -1703.142599999999
$ws.Range("N136").Value = -29758.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1767.6428
$ws.Range("I109").Value = 1437.8334
$ws.Range("J109").Value = 2015
$ws.Range("K109").Value = 4313.5002
$ws.Range("L109").Value = 6045
$ws.Range("M109").Value = -3273.5002
$ws.Range("N109").Value = -8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 115212.22
$ws.Range("I132").Value = 2333.3333
$ws.Range("J132").Value = 171651.67
$ws.Range("K132").Value = 6999.999899999999
$ws.Range("L132").Value = 514955.01
$ws.Range("M132").Value = -4469.999899999999
$ws.Range("N132").Value = -520015.01

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 520103.5
$ws.Range("J69").Value = 520103.5
$ws.Range("L69").Value = 520103.5
$ws.Range("N69").Value = -521725.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H72").Value = 520103.5
$ws.Range("J72").Value = 520103.5
$ws.Range("L72").Value = 1560310.5
$ws.Range("N72").Value = -1568422.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 49615.547
$ws.Range("I132").Value = 79977.234
$ws.Range("J132").Value = 5759.778
$ws.Range("K132").Value = 239931.702
$ws.Range("L132").Value = 17279.334
$ws.Range("M132").Value = -237401.702
$ws.Range("N132").Value = -22339.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2276.5
$ws.Range("I136").Value = 1681.2142
$ws.Range("J136").Value = 4360
$ws.Range("K136").Value = 5043.642599999999
$ws.Range("L136").Value = 13080
$ws.Range("M136").Value = -2493.642599999999
$ws.Range("N136").Value = -18180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 27516.666
$ws.Range("J92").Value = 27516.666
$ws.Range("L92").Value = 27516.666
$ws.Range("N92").Value = -32508.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1054.2258
$ws.Range("I107").Value = 889.94446
$ws.Range("J107").Value = 1281.6923
$ws.Range("K107").Value = 2669.83338
$ws.Range("L107").Value = 3845.0769
$ws.Range("M107").Value = -749.83338
$ws.Range("N107").Value = -7685.0769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1560.2572
$ws.Range("I132").Value = 1295.9333
$ws.Range("J132").Value = 3146.2
$ws.Range("K132").Value = 3887.7999
$ws.Range("L132").Value = 9438.599999999999
$ws.Range("M132").Value = -1357.7999
$ws.Range("N132").Value = -14498.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3878722.5
$ws.Range("I136").Value = 4465654
$ws.Range("J136").Value = 2000542
$ws.Range("K136").Value = 13396962
$ws.Range("L136").Value = 2000542
$ws.Range("M136").Value = -13394412
$ws.Range("N136").Value = -6006726
